$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 728.2857
$ws.Range("J19").Value = 475
$ws.Range("L19").Value = 475
$ws.Range("N19").Value = -825
$ws.Range("H76").Value = 4720
$ws.Range("I76").Value = 4883.1665
$ws.Range("K76").Value = 4883.1665
$ws.Range("M76").Value = -4568.1665
$ws.Range("H79").Value = 4720
$ws.Range("I79").Value = 4883.1665
$ws.Range("K79").Value = 4883.1665
$ws.Range("M79").Value = -3791.1665
$ws.Range("H92").Value = 1042.5
$ws.Range("I92").Value = 1112.6957
$ws.Range("K92").Value = 1112.6957
$ws.Range("M92").Value = 135.3043
$ws.Range("H98").Value = 6004.75
$ws.Range("I98").Value = 5903.4116
$ws.Range("J98").Value = 6579
$ws.Range("K98").Value = 5903.4116
$ws.Range("L98").Value = 6579
$ws.Range("M98").Value = -4405.4116
$ws.Range("N98").Value = -9575
$ws.Range("H116").Value = 2294.3684
$ws.Range("I116").Value = 1892
$ws.Range("J116").Value = 2847.625
$ws.Range("K116").Value = 1892
$ws.Range("L116").Value = 2847.625
$ws.Range("M116").Value = 1550
$ws.Range("N116").Value = -9731.625
$ws.Range("H122").Value = 6004.75
$ws.Range("I122").Value = 5903.4116
$ws.Range("J122").Value = 6579
$ws.Range("K122").Value = 17710.2348
$ws.Range("L122").Value = 19737
$ws.Range("M122").Value = -15260.2348
$ws.Range("N122").Value = -24637
$ws.Range("H132").Value = 3465.7778
$ws.Range("I132").Value = 3024.875
$ws.Range("K132").Value = 9074.625
$ws.Range("M132").Value = -6544.625
$ws.Range("H137").Value = 10345.958
$ws.Range("I137").Value = 1334.7
$ws.Range("J137").Value = 16782.572
$ws.Range("K137").Value = 4004.1
$ws.Range("L137").Value = 50347.716
$ws.Range("M137").Value = -1454.1
$ws.Range("N137").Value = -55447.716
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 38041.25
$ws.Range("J43").Value = 40722
$ws.Range("L43").Value = 40722
$ws.Range("N43").Value = -41348
$ws.Range("H45").Value = 2240.4614
$ws.Range("I45").Value = 854
$ws.Range("K45").Value = 854
$ws.Range("M45").Value = -477
$ws.Range("H102").Value = 1936.5
$ws.Range("I102").Value = 1879.6
$ws.Range("K102").Value = 1879.6
$ws.Range("M102").Value = -257.5999999999999
$ws.Range("H110").Value = 6902.25
$ws.Range("I110").Value = 2392
$ws.Range("J110").Value = 11412.5
$ws.Range("K110").Value = 2392
$ws.Range("L110").Value = 11412.5
$ws.Range("M110").Value = -347
$ws.Range("N110").Value = -15502.5
$ws.Range("H139").Value = 82374.5
$ws.Range("J139").Value = 82374.5
$ws.Range("L139").Value = 82374.5
$ws.Range("N139").Value = -92654.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 360.125
$ws.Range("I22").Value = 358.13333
$ws.Range("K22").Value = 358.13333
$ws.Range("M22").Value = -185.13333
$ws.Range("H94").Value = 20621.219
$ws.Range("I94").Value = 12161.125
$ws.Range("J94").Value = 46001.5
$ws.Range("K94").Value = 12161.125
$ws.Range("L94").Value = 46001.5
$ws.Range("M94").Value = -11710.125
$ws.Range("N94").Value = -46903.5
$ws.Range("H105").Value = 106603.945
$ws.Range("I105").Value = 1499.0667
$ws.Range("K105").Value = 1499.0667
$ws.Range("M105").Value = 247.9332999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3790262.2
$ws.Range("I31").Value = 9092625
$ws.Range("K31").Value = 9092625
$ws.Range("M31").Value = -9092330
$ws.Range("H34").Value = 3790262.2
$ws.Range("I34").Value = 9092625
$ws.Range("K34").Value = 9092625
$ws.Range("M34").Value = -9092423
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 550.25
$ws.Range("J22").Value = 566.6667
$ws.Range("L22").Value = 1700.0001
$ws.Range("N22").Value = -2038.0001
$ws.Range("H26").Value = 45.166668
$ws.Range("I26").Value = 44.8
$ws.Range("K26").Value = 134.4
$ws.Range("M26").Value = 153.6
$ws.Range("H27").Value = 550.25
$ws.Range("J27").Value = 566.6667
$ws.Range("L27").Value = 1700.0001
$ws.Range("N27").Value = -1904.0001
$ws.Range("H31").Value = 2775.2222
$ws.Range("I31").Value = 2747.125
$ws.Range("K31").Value = 8241.375
$ws.Range("M31").Value = -7953.375
$ws.Range("H34").Value = 2990.1724
$ws.Range("J34").Value = 3347.32
$ws.Range("L34").Value = 10041.96
$ws.Range("N34").Value = -10209.96
$ws.Range("H39").Value = 5488.4
$ws.Range("J39").Value = 5887.1113
$ws.Range("L39").Value = 17661.3339
$ws.Range("N39").Value = -18249.3339
$ws.Range("H51").Value = 2671.4285
$ws.Range("I51").Value = 2425
$ws.Range("J51").Value = 3000
$ws.Range("K51").Value = 7275
$ws.Range("L51").Value = 9000
$ws.Range("M51").Value = -6815
$ws.Range("N51").Value = -9920
$ws.Range("H55").Value = 5277.778
$ws.Range("J55").Value = 5277.778
$ws.Range("L55").Value = 15833.334
$ws.Range("N55").Value = -16187.334
$ws.Range("H56").Value = 5866.3335
$ws.Range("I56").Value = 5866.3335
$ws.Range("K56").Value = 5866.3335
$ws.Range("M56").Value = -5336.3335
$ws.Range("H60").Value = 2316.4285
$ws.Range("I60").Value = 3189.4
$ws.Range("J60").Value = 134
$ws.Range("K60").Value = 9568.200000000001
$ws.Range("L60").Value = 402
$ws.Range("M60").Value = -9317.200000000001
$ws.Range("N60").Value = -904
$ws.Range("H64").Value = 18300.834
$ws.Range("I64").Value = 16961
$ws.Range("K64").Value = 50883
$ws.Range("M64").Value = -50613
$ws.Range("H67").Value = 18300.834
$ws.Range("I67").Value = 16961
$ws.Range("K67").Value = 50883
$ws.Range("M67").Value = -49947
$ws.Range("H131").Value = 1665.08
$ws.Range("I131").Value = 1429.75
$ws.Range("J131").Value = 1882.3077
$ws.Range("K131").Value = 4289.25
$ws.Range("L131").Value = 5646.9231
$ws.Range("M131").Value = 750.75
$ws.Range("N131").Value = -15726.9231
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6741.5454
$ws.Range("I70").Value = 5750
$ws.Range("J70").Value = 6961.8887
$ws.Range("K70").Value = 5750
$ws.Range("L70").Value = 6961.8887
$ws.Range("M70").Value = -5480
$ws.Range("N70").Value = -7501.8887
$ws.Range("H73").Value = 6741.5454
$ws.Range("I73").Value = 5750
$ws.Range("J73").Value = 6961.8887
$ws.Range("K73").Value = 5750
$ws.Range("L73").Value = 6961.8887
$ws.Range("M73").Value = -4814
$ws.Range("N73").Value = -8833.8887
$ws.Range("H80").Value = 2496.875
$ws.Range("I80").Value = 2329
$ws.Range("K80").Value = 2329
$ws.Range("M80").Value = -1331
$ws.Range("H83").Value = 2496.875
$ws.Range("I83").Value = 2329
$ws.Range("K83").Value = 11645
$ws.Range("M83").Value = -6653
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 14802.333
$ws.Range("I82").Value = 3342.4443
$ws.Range("J82").Value = 83561.664
$ws.Range("K82").Value = 3342.4443
$ws.Range("L82").Value = 83561.664
$ws.Range("M82").Value = -2981.4443
$ws.Range("N82").Value = -84283.664
$ws.Range("H85").Value = 14802.333
$ws.Range("I85").Value = 3342.4443
$ws.Range("J85").Value = 83561.664
$ws.Range("K85").Value = 3342.4443
$ws.Range("L85").Value = 83561.664
$ws.Range("M85").Value = -2094.4443
$ws.Range("N85").Value = -86057.664
$ws.Range("H93").Value = 2825.1765
$ws.Range("I93").Value = 2060.3914
$ws.Range("J93").Value = 4424.273
$ws.Range("K93").Value = 2060.3914
$ws.Range("L93").Value = 4424.273
$ws.Range("M93").Value = -812.3914
$ws.Range("N93").Value = -6920.273
$ws.Range("H122").Value = 34048.2
$ws.Range("I122").Value = 21925.7
$ws.Range("K122").Value = 65777.10000000001
$ws.Range("M122").Value = -63327.10000000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 16200
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H52").Value = 24875
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H62").Value = 4464.2856
$ws.Range("H65").Value = 4464.2856
